$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.943.55'
$ws.Range('E2').Value = '  -0.68%  '

$ws.Range('D3').Value = '2.354.32'
$ws.Range('E3').Value = '  -0.63%  '

$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('B5').Value = 'XRP'
$ws.Range('C5').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.675'
$ws.Range('D5').Style = $ws.Range('B5').Style
$ws.Range('E5').Value = '  -2.59%  '

$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '240.48'
$ws.Range('D6').Style = $ws.Range('B6').Style
$ws.Range('E6').Value = '  -1.30%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '73.49'
$ws.Range('D7').Style = $ws.Range('B7').Style
$ws.Range('E7').Value = '  -1.20%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.604'
$ws.Range('D9').Style = $ws.Range('B9').Style
$ws.Range('E9').Value = '  +0.12%  '

$ws.Range('E10').Value = '  -2.84%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '58.46'
$ws.Range('D11').Style = $ws.Range('B11').Style
$ws.Range('E11').Value = '  +0.71%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '34.05'
$ws.Range('D12').Style = $ws.Range('B12').Style
$ws.Range('E12').Value = '  +6.35%  '

$ws.Range('E13').Value = '  -2.27%  '

$ws.Range('E14').Value = '  -0.08%  '

$ws.Range('D15').Value = '2.704.95'
$ws.Range('E15').Value = '  -0.60%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '16.38'
$ws.Range('D16').Style = $ws.Range('B16').Style
$ws.Range('E16').Value = '  -4.47%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.913'
$ws.Range('D17').Style = $ws.Range('B17').Style
$ws.Range('E17').Value = '  -0.95%  '

$ws.Range('D18').Value = '2.351.35'
$ws.Range('E18').Value = '  -0.60%  '

$ws.Range('D19').Value = '43.858.13'
$ws.Range('E19').Value = '  -0.89%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0000103'
$ws.Range('D20').Style = $ws.Range('B20').Style
$ws.Range('E20').Value = '  -1.60%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.72'
$ws.Range('D21').Style = $ws.Range('B21').Style
$ws.Range('E21').Value = '  -0.51%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '77.83'
$ws.Range('D22').Style = $ws.Range('B22').Style
$ws.Range('E22').Value = '  -1.37%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '257.45'
$ws.Range('D23').Style = $ws.Range('B23').Style
$ws.Range('E23').Value = '  -0.13%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.93'
$ws.Range('D24').Style = $ws.Range('B24').Style
$ws.Range('E24').Value = '  +15.42%  '

$ws.Range('B25').Value = 'WEMIXToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.80'
$ws.Range('D25').Style = $ws.Range('B25').Style
$ws.Range('E25').Value = '  +1.21%  '

$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = $ws.Range('B26').Style
$ws.Range('E26').Value = '  -0.04%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.50'
$ws.Range('D27').Style = $ws.Range('B27').Style
$ws.Range('E27').Value = '  -2.65%  '

$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.67'
$ws.Range('D28').Style = $ws.Range('B28').Style
$ws.Range('E28').Value = '  -2.00%  '

$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.38'
$ws.Range('D29').Style = $ws.Range('B29').Style
$ws.Range('E29').Value = '  +2.72%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '22.71'
$ws.Range('D30').Style = $ws.Range('B30').Style
$ws.Range('E30').Value = '  -0.28%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '177.86'
$ws.Range('D31').Style = $ws.Range('B31').Style
$ws.Range('E31').Value = '  +1.81%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.130'
$ws.Range('D32').Style = $ws.Range('B32').Style
$ws.Range('E32').Value = '  -0.82%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.136'
$ws.Range('D33').Style = $ws.Range('B33').Style
$ws.Range('E33').Value = '  -0.38%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0758'
$ws.Range('D34').Style = $ws.Range('B34').Style
$ws.Range('E34').Value = '  -0.49%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.47'
$ws.Range('D36').Style = $ws.Range('B36').Style
$ws.Range('E36').Value = '  +0.28%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.83'
$ws.Range('D37').Style = $ws.Range('B37').Style
$ws.Range('E37').Value = '  -2.17%  '

$ws.Range('E38').Value = '  -3.41%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.42'
$ws.Range('D39').Style = $ws.Range('B39').Style
$ws.Range('E39').Value = '  -2.49%  '

$ws.Range('E40').Value = '  +0.93%  '

$ws.Range('B41').Value = 'FTXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.34'
$ws.Range('D41').Style = $ws.Range('B41').Style
$ws.Range('E41').Value = '  +19.53%  '

$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '66.74'
$ws.Range('D42').Style = $ws.Range('B42').Style
$ws.Range('E42').Value = '  +24.75%  '

$ws.Range('E43').Value = '  +10.75%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '9.12'
$ws.Range('D44').Style = $ws.Range('B44').Style
$ws.Range('E44').Value = '  -0.03%  '

$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.204'
$ws.Range('D45').Style = $ws.Range('B45').Style
$ws.Range('E45').Value = '  +1.75%  '

$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '19.06'
$ws.Range('D46').Style = $ws.Range('B46').Style
$ws.Range('E46').Value = '  -0.61%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.52'
$ws.Range('D47').Style = $ws.Range('B47').Style
$ws.Range('E47').Value = '  +0.46%  '

$ws.Range('E48').Value = '  -0.37%  '

$ws.Range('E50').Value = '  -2.57%  '

$ws.Range('E51').Value = '  -2.23%  '
